$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $c = $ws.Range($addr)
    $c.Formula = '="' + $text + '"'
    $c.Copy()
    $c.PasteSpecial(-4163)
}

Set-TextValue "B11" "32.92"
Set-TextValue "C11" "2.88"
Set-TextValue "B12" "28.62"
Set-TextValue "C12" "40.08"
Set-TextValue "C36" "8.03"
Set-TextValue "D36" "99.72"

$excel.CutCopyMode = 0
